$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$cell = $table.Cell(1, 1)
$rng = $cell.Range
$rng.Find.Execute("38+41=", $true, $false, $false, $false, $false, $true, 1, $false, "47-25=", 2) | Out-Null
$cell = $table.Cell(1, 2)
$rng = $cell.Range
$rng.Find.Execute("7+80=", $true, $false, $false, $false, $false, $true, 1, $false, "62+22=", 2) | Out-Null
$cell = $table.Cell(1, 3)
$rng = $cell.Range
$rng.Find.Execute("31-4=", $true, $false, $false, $false, $false, $true, 1, $false, "48+29=", 2) | Out-Null
$cell = $table.Cell(1, 4)
$rng = $cell.Range
$rng.Find.Execute("1+8=", $true, $false, $false, $false, $false, $true, 1, $false, "42-36=", 2) | Out-Null
$cell = $table.Cell(1, 5)
$rng = $cell.Range
$rng.Find.Execute("23+1=", $true, $false, $false, $false, $false, $true, 1, $false, "58-25=", 2) | Out-Null
$cell = $table.Cell(2, 1)
$rng = $cell.Range
$rng.Find.Execute("91-62=", $true, $false, $false, $false, $false, $true, 1, $false, "49+17=", 2) | Out-Null
$cell = $table.Cell(2, 2)
$rng = $cell.Range
$rng.Find.Execute("51+14=", $true, $false, $false, $false, $false, $true, 1, $false, "69-67=", 2) | Out-Null
$cell = $table.Cell(2, 3)
$rng = $cell.Range
$rng.Find.Execute("78-44=", $true, $false, $false, $false, $false, $true, 1, $false, "78+20=", 2) | Out-Null
$cell = $table.Cell(2, 4)
$rng = $cell.Range
$rng.Find.Execute("24-10=", $true, $false, $false, $false, $false, $true, 1, $false, "23+22=", 2) | Out-Null
$cell = $table.Cell(2, 5)
$rng = $cell.Range
$rng.Find.Execute("36+60=", $true, $false, $false, $false, $false, $true, 1, $false, "35-15=", 2) | Out-Null
$cell = $table.Cell(3, 1)
$rng = $cell.Range
$rng.Find.Execute("8+79=", $true, $false, $false, $false, $false, $true, 1, $false, "85-17=", 2) | Out-Null
$cell = $table.Cell(3, 2)
$rng = $cell.Range
$rng.Find.Execute("2+90=", $true, $false, $false, $false, $false, $true, 1, $false, "17-11=", 2) | Out-Null
$cell = $table.Cell(3, 3)
$rng = $cell.Range
$rng.Find.Execute("81-67=", $true, $false, $false, $false, $false, $true, 1, $false, "81-17=", 2) | Out-Null
$cell = $table.Cell(3, 4)
$rng = $cell.Range
$rng.Find.Execute("40+19=", $true, $false, $false, $false, $false, $true, 1, $false, "35+2=", 2) | Out-Null
$cell = $table.Cell(3, 5)
$rng = $cell.Range
$rng.Find.Execute("35-11=", $true, $false, $false, $false, $false, $true, 1, $false, "14+15=", 2) | Out-Null
$cell = $table.Cell(4, 1)
$rng = $cell.Range
$rng.Find.Execute("64-57=", $true, $false, $false, $false, $false, $true, 1, $false, "60-46=", 2) | Out-Null
$cell = $table.Cell(4, 2)
$rng = $cell.Range
$rng.Find.Execute("66-1=", $true, $false, $false, $false, $false, $true, 1, $false, "91-48=", 2) | Out-Null
$cell = $table.Cell(4, 3)
$rng = $cell.Range
$rng.Find.Execute("90-50=", $true, $false, $false, $false, $false, $true, 1, $false, "84-68=", 2) | Out-Null
$cell = $table.Cell(4, 4)
$rng = $cell.Range
$rng.Find.Execute("92-2=", $true, $false, $false, $false, $false, $true, 1, $false, "18-17=", 2) | Out-Null
$cell = $table.Cell(4, 5)
$rng = $cell.Range
$rng.Find.Execute("47-5=", $true, $false, $false, $false, $false, $true, 1, $false, "44-0=", 2) | Out-Null
$cell = $table.Cell(5, 1)
$rng = $cell.Range
$rng.Find.Execute("84-22=", $true, $false, $false, $false, $false, $true, 1, $false, "32-20=", 2) | Out-Null
$cell = $table.Cell(5, 2)
$rng = $cell.Range
$rng.Find.Execute("87-76=", $true, $false, $false, $false, $false, $true, 1, $false, "24+71=", 2) | Out-Null
$cell = $table.Cell(5, 3)
$rng = $cell.Range
$rng.Find.Execute("71-25=", $true, $false, $false, $false, $false, $true, 1, $false, "80-33=", 2) | Out-Null
$cell = $table.Cell(5, 4)
$rng = $cell.Range
$rng.Find.Execute("81-8=", $true, $false, $false, $false, $false, $true, 1, $false, "85-71=", 2) | Out-Null
$cell = $table.Cell(5, 5)
$rng = $cell.Range
$rng.Find.Execute("77-47=", $true, $false, $false, $false, $false, $true, 1, $false, "75-2=", 2) | Out-Null
$cell = $table.Cell(6, 1)
$rng = $cell.Range
$rng.Find.Execute("16+32=", $true, $false, $false, $false, $false, $true, 1, $false, "34+65=", 2) | Out-Null
$cell = $table.Cell(6, 2)
$rng = $cell.Range
$rng.Find.Execute("44+15=", $true, $false, $false, $false, $false, $true, 1, $false, "89-2=", 2) | Out-Null
$cell = $table.Cell(6, 3)
$rng = $cell.Range
$rng.Find.Execute("91-28=", $true, $false, $false, $false, $false, $true, 1, $false, "67+9=", 2) | Out-Null
$cell = $table.Cell(6, 4)
$rng = $cell.Range
$rng.Find.Execute("99-76=", $true, $false, $false, $false, $false, $true, 1, $false, "16+19=", 2) | Out-Null
$cell = $table.Cell(6, 5)
$rng = $cell.Range
$rng.Find.Execute("52+8=", $true, $false, $false, $false, $false, $true, 1, $false, "73+0=", 2) | Out-Null
$cell = $table.Cell(7, 1)
$rng = $cell.Range
$rng.Find.Execute("39-5=", $true, $false, $false, $false, $false, $true, 1, $false, "81-69=", 2) | Out-Null
$cell = $table.Cell(7, 2)
$rng = $cell.Range
$rng.Find.Execute("4+55=", $true, $false, $false, $false, $false, $true, 1, $false, "62-58=", 2) | Out-Null
$cell = $table.Cell(7, 3)
$rng = $cell.Range
$rng.Find.Execute("49-13=", $true, $false, $false, $false, $false, $true, 1, $false, "90-8=", 2) | Out-Null
$cell = $table.Cell(7, 4)
$rng = $cell.Range
$rng.Find.Execute("27+69=", $true, $false, $false, $false, $false, $true, 1, $false, "84-56=", 2) | Out-Null
$cell = $table.Cell(7, 5)
$rng = $cell.Range
$rng.Find.Execute("35+22=", $true, $false, $false, $false, $false, $true, 1, $false, "84+10=", 2) | Out-Null
$cell = $table.Cell(8, 1)
$rng = $cell.Range
$rng.Find.Execute("60+26=", $true, $false, $false, $false, $false, $true, 1, $false, "67+1=", 2) | Out-Null
$cell = $table.Cell(8, 2)
$rng = $cell.Range
$rng.Find.Execute("73-43=", $true, $false, $false, $false, $false, $true, 1, $false, "12+1=", 2) | Out-Null
$cell = $table.Cell(8, 3)
$rng = $cell.Range
$rng.Find.Execute("63-27=", $true, $false, $false, $false, $false, $true, 1, $false, "24-2=", 2) | Out-Null
$cell = $table.Cell(8, 4)
$rng = $cell.Range
$rng.Find.Execute("69-66=", $true, $false, $false, $false, $false, $true, 1, $false, "56-36=", 2) | Out-Null
$cell = $table.Cell(8, 5)
$rng = $cell.Range
$rng.Find.Execute("42-5=", $true, $false, $false, $false, $false, $true, 1, $false, "62-55=", 2) | Out-Null
$cell = $table.Cell(9, 1)
$rng = $cell.Range
$rng.Find.Execute("21+77=", $true, $false, $false, $false, $false, $true, 1, $false, "71+10=", 2) | Out-Null
$cell = $table.Cell(9, 2)
$rng = $cell.Range
$rng.Find.Execute("8+54=", $true, $false, $false, $false, $false, $true, 1, $false, "13+40=", 2) | Out-Null
$cell = $table.Cell(9, 3)
$rng = $cell.Range
$rng.Find.Execute("50-42=", $true, $false, $false, $false, $false, $true, 1, $false, "9+39=", 2) | Out-Null
$cell = $table.Cell(9, 4)
$rng = $cell.Range
$rng.Find.Execute("19+78=", $true, $false, $false, $false, $false, $true, 1, $false, "87-85=", 2) | Out-Null
$cell = $table.Cell(9, 5)
$rng = $cell.Range
$rng.Find.Execute("38-35=", $true, $false, $false, $false, $false, $true, 1, $false, "91-17=", 2) | Out-Null
$cell = $table.Cell(10, 1)
$rng = $cell.Range
$rng.Find.Execute("92-34=", $true, $false, $false, $false, $false, $true, 1, $false, "87-32=", 2) | Out-Null
$cell = $table.Cell(10, 2)
$rng = $cell.Range
$rng.Find.Execute("20-19=", $true, $false, $false, $false, $false, $true, 1, $false, "19+15=", 2) | Out-Null
$cell = $table.Cell(10, 3)
$rng = $cell.Range
$rng.Find.Execute("33+6=", $true, $false, $false, $false, $false, $true, 1, $false, "29+4=", 2) | Out-Null
$cell = $table.Cell(10, 4)
$rng = $cell.Range
$rng.Find.Execute("54+13=", $true, $false, $false, $false, $false, $true, 1, $false, "3+37=", 2) | Out-Null
$cell = $table.Cell(10, 5)
$rng = $cell.Range
$rng.Find.Execute("66-32=", $true, $false, $false, $false, $false, $true, 1, $false, "11+61=", 2) | Out-Null
$cell = $table.Cell(11, 1)
$rng = $cell.Range
$rng.Find.Execute("35+28=", $true, $false, $false, $false, $false, $true, 1, $false, "91-2=", 2) | Out-Null
$cell = $table.Cell(11, 2)
$rng = $cell.Range
$rng.Find.Execute("79-68=", $true, $false, $false, $false, $false, $true, 1, $false, "88-61=", 2) | Out-Null
$cell = $table.Cell(11, 3)
$rng = $cell.Range
$rng.Find.Execute("82-2=", $true, $false, $false, $false, $false, $true, 1, $false, "80+16=", 2) | Out-Null
$cell = $table.Cell(11, 4)
$rng = $cell.Range
$rng.Find.Execute("4+84=", $true, $false, $false, $false, $false, $true, 1, $false, "5-4=", 2) | Out-Null
$cell = $table.Cell(11, 5)
$rng = $cell.Range
$rng.Find.Execute("21+40=", $true, $false, $false, $false, $false, $true, 1, $false, "20-10=", 2) | Out-Null
$cell = $table.Cell(12, 1)
$rng = $cell.Range
$rng.Find.Execute("8+13=", $true, $false, $false, $false, $false, $true, 1, $false, "68-51=", 2) | Out-Null
$cell = $table.Cell(12, 2)
$rng = $cell.Range
$rng.Find.Execute("25+65=", $true, $false, $false, $false, $false, $true, 1, $false, "24+63=", 2) | Out-Null
$cell = $table.Cell(12, 3)
$rng = $cell.Range
$rng.Find.Execute("21-7=", $true, $false, $false, $false, $false, $true, 1, $false, "58+33=", 2) | Out-Null
$cell = $table.Cell(12, 4)
$rng = $cell.Range
$rng.Find.Execute("59+26=", $true, $false, $false, $false, $false, $true, 1, $false, "70-28=", 2) | Out-Null
$cell = $table.Cell(12, 5)
$rng = $cell.Range
$rng.Find.Execute("40+52=", $true, $false, $false, $false, $false, $true, 1, $false, "49+31=", 2) | Out-Null
$cell = $table.Cell(13, 1)
$rng = $cell.Range
$rng.Find.Execute("47+13=", $true, $false, $false, $false, $false, $true, 1, $false, "98-67=", 2) | Out-Null
$cell = $table.Cell(13, 2)
$rng = $cell.Range
$rng.Find.Execute("23+36=", $true, $false, $false, $false, $false, $true, 1, $false, "93-1=", 2) | Out-Null
$cell = $table.Cell(13, 3)
$rng = $cell.Range
$rng.Find.Execute("50+39=", $true, $false, $false, $false, $false, $true, 1, $false, "15+72=", 2) | Out-Null
$cell = $table.Cell(13, 4)
$rng = $cell.Range
$rng.Find.Execute("95-55=", $true, $false, $false, $false, $false, $true, 1, $false, "4+37=", 2) | Out-Null
$cell = $table.Cell(13, 5)
$rng = $cell.Range
$rng.Find.Execute("48+36=", $true, $false, $false, $false, $false, $true, 1, $false, "89-45=", 2) | Out-Null
$cell = $table.Cell(14, 1)
$rng = $cell.Range
$rng.Find.Execute("2+63=", $true, $false, $false, $false, $false, $true, 1, $false, "39-34=", 2) | Out-Null
$cell = $table.Cell(14, 2)
$rng = $cell.Range
$rng.Find.Execute("21+71=", $true, $false, $false, $false, $false, $true, 1, $false, "78+13=", 2) | Out-Null
$cell = $table.Cell(14, 3)
$rng = $cell.Range
$rng.Find.Execute("17+40=", $true, $false, $false, $false, $false, $true, 1, $false, "75-15=", 2) | Out-Null
$cell = $table.Cell(14, 4)
$rng = $cell.Range
$rng.Find.Execute("80+14=", $true, $false, $false, $false, $false, $true, 1, $false, "33+9=", 2) | Out-Null
$cell = $table.Cell(14, 5)
$rng = $cell.Range
$rng.Find.Execute("19+60=", $true, $false, $false, $false, $false, $true, 1, $false, "21+20=", 2) | Out-Null
$cell = $table.Cell(15, 1)
$rng = $cell.Range
$rng.Find.Execute("6+60=", $true, $false, $false, $false, $false, $true, 1, $false, "20+23=", 2) | Out-Null
$cell = $table.Cell(15, 2)
$rng = $cell.Range
$rng.Find.Execute("65-59=", $true, $false, $false, $false, $false, $true, 1, $false, "72-12=", 2) | Out-Null
$cell = $table.Cell(15, 3)
$rng = $cell.Range
$rng.Find.Execute("3+17=", $true, $false, $false, $false, $false, $true, 1, $false, "63-6=", 2) | Out-Null
$cell = $table.Cell(15, 4)
$rng = $cell.Range
$rng.Find.Execute("40-40=", $true, $false, $false, $false, $false, $true, 1, $false, "93+1=", 2) | Out-Null
$cell = $table.Cell(15, 5)
$rng = $cell.Range
$rng.Find.Execute("79+8=", $true, $false, $false, $false, $false, $true, 1, $false, "68-61=", 2) | Out-Null
$cell = $table.Cell(16, 1)
$rng = $cell.Range
$rng.Find.Execute("25+36=", $true, $false, $false, $false, $false, $true, 1, $false, "55+2=", 2) | Out-Null
$cell = $table.Cell(16, 2)
$rng = $cell.Range
$rng.Find.Execute("13+37=", $true, $false, $false, $false, $false, $true, 1, $false, "94-91=", 2) | Out-Null
$cell = $table.Cell(16, 3)
$rng = $cell.Range
$rng.Find.Execute("29-10=", $true, $false, $false, $false, $false, $true, 1, $false, "5+36=", 2) | Out-Null
$cell = $table.Cell(16, 4)
$rng = $cell.Range
$rng.Find.Execute("90-61=", $true, $false, $false, $false, $false, $true, 1, $false, "77-54=", 2) | Out-Null
$cell = $table.Cell(16, 5)
$rng = $cell.Range
$rng.Find.Execute("26+61=", $true, $false, $false, $false, $false, $true, 1, $false, "60-51=", 2) | Out-Null
$cell = $table.Cell(17, 1)
$rng = $cell.Range
$rng.Find.Execute("36+24=", $true, $false, $false, $false, $false, $true, 1, $false, "58-28=", 2) | Out-Null
$cell = $table.Cell(17, 2)
$rng = $cell.Range
$rng.Find.Execute("62+30=", $true, $false, $false, $false, $false, $true, 1, $false, "78-33=", 2) | Out-Null
$cell = $table.Cell(17, 3)
$rng = $cell.Range
$rng.Find.Execute("34+48=", $true, $false, $false, $false, $false, $true, 1, $false, "18+31=", 2) | Out-Null
$cell = $table.Cell(17, 4)
$rng = $cell.Range
$rng.Find.Execute("95-38=", $true, $false, $false, $false, $false, $true, 1, $false, "94-72=", 2) | Out-Null
$cell = $table.Cell(17, 5)
$rng = $cell.Range
$rng.Find.Execute("63-15=", $true, $false, $false, $false, $false, $true, 1, $false, "61-29=", 2) | Out-Null
$cell = $table.Cell(18, 1)
$rng = $cell.Range
$rng.Find.Execute("58-0=", $true, $false, $false, $false, $false, $true, 1, $false, "52-41=", 2) | Out-Null
$cell = $table.Cell(18, 2)
$rng = $cell.Range
$rng.Find.Execute("78+19=", $true, $false, $false, $false, $false, $true, 1, $false, "89-7=", 2) | Out-Null
$cell = $table.Cell(18, 3)
$rng = $cell.Range
$rng.Find.Execute("87-0=", $true, $false, $false, $false, $false, $true, 1, $false, "34-2=", 2) | Out-Null
$cell = $table.Cell(18, 4)
$rng = $cell.Range
$rng.Find.Execute("97-94=", $true, $false, $false, $false, $false, $true, 1, $false, "2+48=", 2) | Out-Null
$cell = $table.Cell(18, 5)
$rng = $cell.Range
$rng.Find.Execute("93-85=", $true, $false, $false, $false, $false, $true, 1, $false, "21-4=", 2) | Out-Null
$cell = $table.Cell(19, 1)
$rng = $cell.Range
$rng.Find.Execute("56+37=", $true, $false, $false, $false, $false, $true, 1, $false, "37+8=", 2) | Out-Null
$cell = $table.Cell(19, 2)
$rng = $cell.Range
$rng.Find.Execute("87+5=", $true, $false, $false, $false, $false, $true, 1, $false, "85-8=", 2) | Out-Null
$cell = $table.Cell(19, 3)
$rng = $cell.Range
$rng.Find.Execute("17+6=", $true, $false, $false, $false, $false, $true, 1, $false, "79-50=", 2) | Out-Null
$cell = $table.Cell(19, 4)
$rng = $cell.Range
$rng.Find.Execute("67+13=", $true, $false, $false, $false, $false, $true, 1, $false, "54-1=", 2) | Out-Null
$cell = $table.Cell(19, 5)
$rng = $cell.Range
$rng.Find.Execute("29+29=", $true, $false, $false, $false, $false, $true, 1, $false, "88-71=", 2) | Out-Null
$cell = $table.Cell(20, 1)
$rng = $cell.Range
$rng.Find.Execute("36-32=", $true, $false, $false, $false, $false, $true, 1, $false, "69-42=", 2) | Out-Null
$cell = $table.Cell(20, 2)
$rng = $cell.Range
$rng.Find.Execute("57-36=", $true, $false, $false, $false, $false, $true, 1, $false, "16+47=", 2) | Out-Null
$cell = $table.Cell(20, 3)
$rng = $cell.Range
$rng.Find.Execute("56-29=", $true, $false, $false, $false, $false, $true, 1, $false, "9-8=", 2) | Out-Null
$cell = $table.Cell(20, 4)
$rng = $cell.Range
$rng.Find.Execute("25+38=", $true, $false, $false, $false, $false, $true, 1, $false, "6+61=", 2) | Out-Null
$cell = $table.Cell(20, 5)
$rng = $cell.Range
$rng.Find.Execute("36-14=", $true, $false, $false, $false, $false, $true, 1, $false, "41+30=", 2) | Out-Null
